# Devis_societe_1 / "Citation" sheet — quarterly quantity correction + Siret fix.
#
# 1. The first line item's quantity (Tableau_Articles row 1, cell B17) is bumped
#    from 32 to 34; this is a table-bound cell so Excel automatically ripples the
#    row's "Montant" formula (F17), the table total (F23) and the grand total
#    (F27) — no need to touch those cells directly, the workbook recalculates on
#    its own (Calculation = xlAutomatic).
# 2. The company's Siret number (just below the VAT/total block) had a typo and
#    gets corrected.
# 3. The user's cursor ends up sitting on the quantity cell they just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the Siret number (shared string shown in C26, next to the
#     "Siret : " label in B26) ---
$ws.Range("C26").Value = "987-654-321 1234"

# --- Bump the quantity on the first article line of Tableau_Articles ---
$ws.Range("B17").Value = 34

# --- Leave the selection on the cell that was just edited ---
$ws.Range("B17").Select()
